# Update the dSF column (column F) with repulled / recalculated data.
# Mirrors the commit "repull data, push all data, mean calculation":
# column F values that used to duplicate column E (dS0) are replaced with the
# freshly pulled final-stat delta values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -2
    4  = 1
    5  = 3
    7  = 2
    9  = 1
    10 = -3
    11 = 1
    12 = -2
    13 = -6
    14 = -2
    15 = 12
    16 = -7
    18 = 2
    19 = -4
    20 = 6
    21 = 1
    23 = 2
    25 = -1
    26 = -1
    27 = -1
    28 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
